$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2560.4
$ws.Range("I31").Value = 2560.4
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 7681.200000000001
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -7451.200000000001
$ws.Range("H43").Value = 905
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 905
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = $null
$ws.Range("M43").Value = 905
$ws.Range("N43").Value = -1043
$ws.Range("H55").Value = 230.66667
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 236.8
$ws.Range("K55").Value = 200
$ws.Range("L55").Value = 236.8
$ws.Range("M55").Value = 14
$ws.Range("N55").Value = -664.8
$ws.Range("H70").Value = 1221.8572
$ws.Range("I70").Value = 513
$ws.Range("J70").Value = 1505.4
$ws.Range("K70").Value = 1539
$ws.Range("L70").Value = 4516.200000000001
$ws.Range("M70").Value = -1269
$ws.Range("N70").Value = -5056.200000000001
$ws.Range("H73").Value = 1221.8572
$ws.Range("I73").Value = 513
$ws.Range("J73").Value = 1505.4
$ws.Range("K73").Value = 1539
$ws.Range("L73").Value = 4516.200000000001
$ws.Range("M73").Value = -603
$ws.Range("N73").Value = -6388.200000000001
$ws.Range("H74").Value = 4725
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 4725
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = $null
$ws.Range("M74").Value = 4725
$ws.Range("N74").Value = -6597
$ws.Range("H77").Value = 4725
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 4725
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = $null
$ws.Range("M77").Value = 23625
$ws.Range("N77").Value = -32985
$ws.Range("H80").Value = 1886.3334
$ws.Range("I80").Value = 1741.6666
$ws.Range("J80").Value = 1958.6666
$ws.Range("K80").Value = 5224.9998
$ws.Range("L80").Value = 5875.9998
$ws.Range("M80").Value = -4226.9998
$ws.Range("N80").Value = -7871.9998
$ws.Range("H83").Value = 1886.3334
$ws.Range("I83").Value = 1741.6666
$ws.Range("J83").Value = 1958.6666
$ws.Range("K83").Value = 15674.9994
$ws.Range("L83").Value = 17627.9994
$ws.Range("M83").Value = -10682.9994
$ws.Range("N83").Value = -27611.9994
$ws.Range("H92").Value = 926.55554
$ws.Range("I92").Value = 917.375
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 917.375
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 330.625
$ws.Range("N92").Value = -3496
$ws.Range("H100").Value = 1999
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null
$ws.Range("H138").Value = 2283.7083
$ws.Range("I138").Value = 851.25
$ws.Range("K138").Value = 2553.75
$ws.Range("M138").Value = 2586.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1137.3334
$ws.Range("I61").Value = 1137.3334
$ws.Range("K61").Value = 1137.3334
$ws.Range("M61").Value = -925.3334
$ws.Range("H63").Value = 6997.25
$ws.Range("J63").Value = 9635
$ws.Range("L63").Value = 9635
$ws.Range("N63").Value = -11007
$ws.Range("H66").Value = 6997.25
$ws.Range("J66").Value = 9635
$ws.Range("L66").Value = 48175
$ws.Range("N66").Value = -55039
$ws.Range("H136").Value = 1137.3334
$ws.Range("I136").Value = 1137.3334
$ws.Range("K136").Value = 3412.0002
$ws.Range("M136").Value = -862.0001999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 148.33333
$ws.Range("I22").Value = 187.5
$ws.Range("K22").Value = 187.5
$ws.Range("M22").Value = 162.5
$ws.Range("H58").Value = 3354.6155
$ws.Range("I58").Value = 2740
$ws.Range("K58").Value = 2740
$ws.Range("M58").Value = -2537
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = $null
$ws.Range("N59").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = $null
$ws.Range("H134").Value = 2368.3157
$ws.Range("I134").Value = 1249.8334
$ws.Range("K134").Value = 3749.5002
$ws.Range("M134").Value = -1214.5002
$ws.Range("H136").Value = 3354.6155
$ws.Range("I136").Value = 2740
$ws.Range("K136").Value = 8220
$ws.Range("M136").Value = -5670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3095.75
$ws.Range("I34").Value = 125
$ws.Range("J34").Value = 3689.9
$ws.Range("K34").Value = 375
$ws.Range("L34").Value = 11069.7
$ws.Range("M34").Value = -291
$ws.Range("N34").Value = -11237.7
$ws.Range("H39").Value = 5972.727
$ws.Range("I39").Value = 7500
$ws.Range("J39").Value = 5633.3335
$ws.Range("K39").Value = 22500
$ws.Range("L39").Value = 16900.0005
$ws.Range("M39").Value = -22206
$ws.Range("N39").Value = -17488.0005
$ws.Range("H55").Value = 2669.9
$ws.Range("J55").Value = 2922.111
$ws.Range("L55").Value = 8766.332999999999
$ws.Range("N55").Value = -9120.332999999999
$ws.Range("H60").Value = 412.25
$ws.Range("I60").Value = 150
$ws.Range("J60").Value = 499.66666
$ws.Range("K60").Value = 450
$ws.Range("L60").Value = 1498.99998
$ws.Range("M60").Value = -199
$ws.Range("N60").Value = -2000.99998
$ws.Range("H139").Value = 1202.8
$ws.Range("J139").Value = 2000
$ws.Range("L139").Value = 6000
$ws.Range("N139").Value = -16280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1705.4445
$ws.Range("I102").Value = 1731.125
$ws.Range("K102").Value = 1731.125
$ws.Range("M102").Value = -109.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2333.3333
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2205
$ws.Range("H27").Value = 2333.3333
$ws.Range("I27").Value = 2500
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = -2393
$ws.Range("H35").Value = 1383.6666
$ws.Range("I35").Value = 575
$ws.Range("K35").Value = 575
$ws.Range("M35").Value = -239
$ws.Range("H40").Value = 8372.799999999999
$ws.Range("I40").Value = 8372.799999999999
$ws.Range("K40").Value = 8372.799999999999
$ws.Range("M40").Value = -8236.799999999999
$ws.Range("H122").Value = 3055.5
$ws.Range("I122").Value = 1111
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 3333
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -883
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 10898.9
$ws.Range("I132").Value = 10898.9
$ws.Range("K132").Value = 32696.7
$ws.Range("M132").Value = -30166.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3122
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
